$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell D1: same label + style as the existing header row (A1:C1)
$ws.Range('C1').Copy() | Out-Null
$ws.Range('D1').PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 4).Value = 'Explanation'

# Explanation text for every data row (plain, unstyled cells like the rest of the data)
$explanations = @{
    2 = 'Matched due to overlap in key terms: java, can, an, looking, assessment, that, minutes, with, in, who'
    3 = 'Matched due to overlap in key terms: java, can, an, looking, assessment, that, minutes, with, in, who'
    4 = 'Matched due to overlap in key terms: java, can, an, looking, assessment, that, minutes, with, in, who'
    5 = 'Matched due to overlap in key terms: insight, looking, science, huge, flexibility, market, ago, candidates, key, best, develop, touch, lead, diversity, workplace, sure, find, excellent, through, business, with, lifeatshl, careeropportunities, than, no, success, years, family, coaching, job, one, is, ability, services, future, around, an, shlhiringtalent, s, solution, transform, talent, are, knowledgeable, description, how, to, become, contribute, experience, us, flexible, careersatshl, fabric, join, workplaces, as, development, encourage, world, inclusion, breaking, something, engineer, me, started, adjustments, possible, career, ongoing, you, your, out, at, about, takes, there, we, goals, inclusivity, employer, shl, when, industry, offer, the, businesses, growth, and, experiences, time, guidance, process, jobopportunities, knowledge, others, shaping, possibilities, from, opportunity, unlock, insights, fun, equity, accessibility, that, power, threads, choose, desirable, more, product, doing, care, technology, get, this, greatest, inspired, continue, products, path, new, people, so, help, by, make, equal, 40, where, culture, scale, be, community, work, offered, assessment, transformational, lives, will, intrinsic, off, making, in, collaboration, acquisition, outcomes, team, package, employee, achieve, application, do, all, better, friendly, benefits, ll, manager, launches, what, on, transformation, ground, a, investment, dei, for, applications, answers, can, part, support, inclusive, essential, inspire, of, range, recruitment, our, diverse, currently, powerhouse'
    6 = 'Matched due to overlap in key terms: insight, looking, science, huge, flexibility, market, ago, candidates, key, best, develop, touch, lead, diversity, workplace, sure, find, excellent, through, business, with, lifeatshl, careeropportunities, than, no, success, years, family, coaching, job, one, is, ability, services, future, around, an, shlhiringtalent, s, solution, transform, talent, are, knowledgeable, description, how, to, become, contribute, experience, us, flexible, careersatshl, fabric, join, workplaces, as, development, encourage, world, inclusion, breaking, something, engineer, me, started, adjustments, possible, career, ongoing, you, your, out, at, about, takes, there, we, goals, inclusivity, employer, shl, when, industry, offer, the, businesses, growth, and, experiences, time, guidance, process, jobopportunities, knowledge, others, shaping, possibilities, from, opportunity, unlock, insights, fun, equity, accessibility, that, power, threads, choose, desirable, more, product, doing, care, technology, get, this, greatest, inspired, continue, products, path, new, people, so, help, by, make, equal, 40, where, culture, scale, be, community, work, offered, assessment, transformational, lives, will, intrinsic, off, making, in, collaboration, acquisition, outcomes, team, package, employee, achieve, application, do, all, better, friendly, benefits, ll, manager, launches, what, on, transformation, ground, a, investment, dei, for, applications, answers, can, part, support, inclusive, essential, inspire, of, range, recruitment, our, diverse, currently, powerhouse'
    7 = 'Matched due to overlap in key terms: insight, looking, science, huge, flexibility, market, ago, candidates, key, best, develop, touch, lead, diversity, workplace, sure, find, excellent, through, business, with, lifeatshl, careeropportunities, than, no, success, years, family, coaching, job, one, is, ability, services, future, around, an, shlhiringtalent, s, solution, transform, talent, are, knowledgeable, description, how, to, become, contribute, experience, us, flexible, careersatshl, fabric, join, workplaces, as, development, encourage, world, inclusion, breaking, something, engineer, me, started, adjustments, possible, career, ongoing, you, your, out, at, about, takes, there, we, goals, inclusivity, employer, shl, when, industry, offer, the, businesses, growth, and, experiences, time, guidance, process, jobopportunities, knowledge, others, shaping, possibilities, from, opportunity, unlock, insights, fun, equity, accessibility, that, power, threads, choose, desirable, more, product, doing, care, technology, get, this, greatest, inspired, continue, products, path, new, people, so, help, by, make, equal, 40, where, culture, scale, be, community, work, offered, assessment, transformational, lives, will, intrinsic, off, making, in, collaboration, acquisition, outcomes, team, package, employee, achieve, application, do, all, better, friendly, benefits, ll, manager, launches, what, on, transformation, ground, a, investment, dei, for, applications, answers, can, part, support, inclusive, essential, inspire, of, range, recruitment, our, diverse, currently, powerhouse'
    8 = 'Matched due to overlap in key terms: what, applications, an, and, are, to, for, tests'
    9 = 'Matched due to overlap in key terms: what, applications, an, and, are, to, for, tests'
    10 = 'Matched due to overlap in key terms: what, applications, an, and, are, to, for, tests'
    11 = 'Matched due to overlap in key terms: insight, science, want, huge, flexibility, market, ago, candidates, key, best, develop, touch, lead, diversity, design, workplace, sure, find, excellent, through, business, with, lifeatshl, than, no, success, years, family, coaching, job, or, one, is, ability, services, around, an, s, solution, transform, talent, are, knowledgeable, description, how, to, become, experience, us, flexible, careersatshl, fabric, join, workplaces, as, development, encourage, world, inclusion, breaking, something, started, getting, adjustments, possible, career, ongoing, you, your, out, at, ready, about, takes, there, we, goals, inclusivity, employer, shl, when, industry, offer, the, businesses, growth, and, required, experiences, time, guidance, process, knowledge, others, possibilities, test, from, opportunity, unlock, quality, insights, fun, equity, accessibility, that, solutions, power, threads, choose, desirable, more, product, care, technology, get, this, greatest, inspired, continue, products, path, new, people, so, help, by, make, equal, functional, 40, where, culture, scale, be, work, offered, transformational, lives, will, intrinsic, off, making, in, collaboration, acquisition, create, team, package, outcomes, employee, achieve, application, do, all, better, benefits, ll, manager, launches, what, on, transformation, ground, a, investment, dei, long, for, applications, answers, can, part, support, inclusive, inspire, creating, strong, of, range, recruitment, our, diverse, currently, powerhouse'
    12 = 'Matched due to overlap in key terms: insight, science, want, huge, flexibility, market, ago, candidates, key, best, develop, touch, lead, diversity, design, workplace, sure, find, excellent, through, business, with, lifeatshl, than, no, success, years, family, coaching, job, or, one, is, ability, services, around, an, s, solution, transform, talent, are, knowledgeable, description, how, to, become, experience, us, flexible, careersatshl, fabric, join, workplaces, as, development, encourage, world, inclusion, breaking, something, started, getting, adjustments, possible, career, ongoing, you, your, out, at, ready, about, takes, there, we, goals, inclusivity, employer, shl, when, industry, offer, the, businesses, growth, and, required, experiences, time, guidance, process, knowledge, others, possibilities, test, from, opportunity, unlock, quality, insights, fun, equity, accessibility, that, solutions, power, threads, choose, desirable, more, product, care, technology, get, this, greatest, inspired, continue, products, path, new, people, so, help, by, make, equal, functional, 40, where, culture, scale, be, work, offered, transformational, lives, will, intrinsic, off, making, in, collaboration, acquisition, create, team, package, outcomes, employee, achieve, application, do, all, better, benefits, ll, manager, launches, what, on, transformation, ground, a, investment, dei, long, for, applications, answers, can, part, support, inclusive, inspire, creating, strong, of, range, recruitment, our, diverse, currently, powerhouse'
    13 = 'Matched due to overlap in key terms: insight, science, want, huge, flexibility, market, ago, candidates, key, best, develop, touch, lead, diversity, design, workplace, sure, find, excellent, through, business, with, lifeatshl, than, no, success, years, family, coaching, job, or, one, is, ability, services, around, an, s, solution, transform, talent, are, knowledgeable, description, how, to, become, experience, us, flexible, careersatshl, fabric, join, workplaces, as, development, encourage, world, inclusion, breaking, something, started, getting, adjustments, possible, career, ongoing, you, your, out, at, ready, about, takes, there, we, goals, inclusivity, employer, shl, when, industry, offer, the, businesses, growth, and, required, experiences, time, guidance, process, knowledge, others, possibilities, test, from, opportunity, unlock, quality, insights, fun, equity, accessibility, that, solutions, power, threads, choose, desirable, more, product, care, technology, get, this, greatest, inspired, continue, products, path, new, people, so, help, by, make, equal, functional, 40, where, culture, scale, be, work, offered, transformational, lives, will, intrinsic, off, making, in, collaboration, acquisition, create, team, package, outcomes, employee, achieve, application, do, all, better, benefits, ll, manager, launches, what, on, transformation, ground, a, investment, dei, long, for, applications, answers, can, part, support, inclusive, inspire, creating, strong, of, range, recruitment, our, diverse, currently, powerhouse'
    14 = 'Matched due to overlap in key terms: an, new, graduates, i, sales, in, for, my'
    15 = 'Matched due to overlap in key terms: an, new, graduates, i, sales, in, for, my'
    16 = 'Matched due to overlap in key terms: an, new, graduates, i, sales, in, for, my'
    17 = 'Matched due to overlap in key terms: other, culture, be, content, best, opportunity, in, e, team, working, driven, that, continuous, with, about, on, we, provide, organization, communication, product, a, members, individuals, one, is, the, for, com, leading, s, and, are, description, of, to, executing, equal, experience, position'
    18 = 'Matched due to overlap in key terms: other, culture, be, content, best, opportunity, in, e, team, working, driven, that, continuous, with, about, on, we, provide, organization, communication, product, a, members, individuals, one, is, the, for, com, leading, s, and, are, description, of, to, executing, equal, experience, position'
    19 = 'Matched due to overlap in key terms: other, culture, be, content, best, opportunity, in, e, team, working, driven, that, continuous, with, about, on, we, provide, organization, communication, product, a, members, individuals, one, is, the, for, com, leading, s, and, are, description, of, to, executing, equal, experience, position'
    20 = 'Matched due to overlap in key terms: and, want, expertise, a, years, i, hire, of, to, with, experience, in'
    21 = 'Matched due to overlap in key terms: and, want, expertise, a, years, i, hire, of, to, with, experience, in'
    22 = 'Matched due to overlap in key terms: and, want, expertise, a, years, i, hire, of, to, with, experience, in'
    23 = 'Matched due to overlap in key terms: insight, looking, science, huge, flexibility, market, ago, candidates, key, best, develop, touch, lead, diversity, workplace, sure, find, excellent, through, business, with, lifeatshl, management, than, no, success, years, family, coaching, job, or, one, is, ability, services, around, an, shlhiringtalent, s, transform, talent, are, knowledgeable, description, how, to, become, experience, below, us, flexible, careersatshl, fabric, join, workplaces, as, development, encourage, world, inclusion, breaking, something, me, started, adjustments, possible, career, ongoing, you, your, out, at, about, takes, there, we, providing, goals, inclusivity, employer, shl, when, industry, offer, the, businesses, growth, and, experiences, time, guidance, process, knowledge, others, possibilities, 1, from, opportunity, unlock, quality, insights, fun, equity, accessibility, that, power, processes, threads, choose, more, product, doing, care, technology, get, this, greatest, inspired, continue, products, path, new, people, so, help, by, make, equal, 40, where, culture, scale, be, work, offered, assessment, transformational, lives, will, intrinsic, off, making, in, collaboration, acquisition, outcomes, team, package, employee, achieve, application, do, all, better, friendly, benefits, ll, manager, launches, what, on, transformation, ground, a, investment, dei, for, applications, answers, can, part, support, inclusive, inspire, strong, of, range, recruitment, our, diverse, currently, powerhouse'
    24 = 'Matched due to overlap in key terms: insight, looking, science, huge, flexibility, market, ago, candidates, key, best, develop, touch, lead, diversity, workplace, sure, find, excellent, through, business, with, lifeatshl, management, than, no, success, years, family, coaching, job, or, one, is, ability, services, around, an, shlhiringtalent, s, transform, talent, are, knowledgeable, description, how, to, become, experience, below, us, flexible, careersatshl, fabric, join, workplaces, as, development, encourage, world, inclusion, breaking, something, me, started, adjustments, possible, career, ongoing, you, your, out, at, about, takes, there, we, providing, goals, inclusivity, employer, shl, when, industry, offer, the, businesses, growth, and, experiences, time, guidance, process, knowledge, others, possibilities, 1, from, opportunity, unlock, quality, insights, fun, equity, accessibility, that, power, processes, threads, choose, more, product, doing, care, technology, get, this, greatest, inspired, continue, products, path, new, people, so, help, by, make, equal, 40, where, culture, scale, be, work, offered, assessment, transformational, lives, will, intrinsic, off, making, in, collaboration, acquisition, outcomes, team, package, employee, achieve, application, do, all, better, friendly, benefits, ll, manager, launches, what, on, transformation, ground, a, investment, dei, for, applications, answers, can, part, support, inclusive, inspire, strong, of, range, recruitment, our, diverse, currently, powerhouse'
    25 = 'Matched due to overlap in key terms: insight, looking, science, huge, flexibility, market, ago, candidates, key, best, develop, touch, lead, diversity, workplace, sure, find, excellent, through, business, with, lifeatshl, management, than, no, success, years, family, coaching, job, or, one, is, ability, services, around, an, shlhiringtalent, s, transform, talent, are, knowledgeable, description, how, to, become, experience, below, us, flexible, careersatshl, fabric, join, workplaces, as, development, encourage, world, inclusion, breaking, something, me, started, adjustments, possible, career, ongoing, you, your, out, at, about, takes, there, we, providing, goals, inclusivity, employer, shl, when, industry, offer, the, businesses, growth, and, experiences, time, guidance, process, knowledge, others, possibilities, 1, from, opportunity, unlock, quality, insights, fun, equity, accessibility, that, power, processes, threads, choose, more, product, doing, care, technology, get, this, greatest, inspired, continue, products, path, new, people, so, help, by, make, equal, 40, where, culture, scale, be, work, offered, assessment, transformational, lives, will, intrinsic, off, making, in, collaboration, acquisition, outcomes, team, package, employee, achieve, application, do, all, better, friendly, benefits, ll, manager, launches, what, on, transformation, ground, a, investment, dei, for, applications, answers, can, part, support, inclusive, inspire, strong, of, range, recruitment, our, diverse, currently, powerhouse'
    26 = 'Matched due to overlap in key terms: manage, us, other, flexible, be, work, skills, deliver, as, not, will, professional, in, team, responsibilities, working, i, that, building, with, benefits, at, management, about, on, issues, providing, we, provide, communication, product, years, more, based, employer, sales, a, technology, or, when, this, is, ability, role, the, for, relevant, services, businesses, an, expert, leading, s, and, organizational, support, internal, opportunities, are, ensuring, teams, 2, by, strong, of, to, equal, experience'
    27 = 'Matched due to overlap in key terms: manage, us, other, flexible, be, work, skills, deliver, as, not, will, professional, in, team, responsibilities, working, i, that, building, with, benefits, at, management, about, on, issues, providing, we, provide, communication, product, years, more, based, employer, sales, a, technology, or, when, this, is, ability, role, the, for, relevant, services, businesses, an, expert, leading, s, and, organizational, support, internal, opportunities, are, ensuring, teams, 2, by, strong, of, to, equal, experience'
    28 = 'Matched due to overlap in key terms: manage, us, other, flexible, be, work, skills, deliver, as, not, will, professional, in, team, responsibilities, working, i, that, building, with, benefits, at, management, about, on, issues, providing, we, provide, communication, product, years, more, based, employer, sales, a, technology, or, when, this, is, ability, role, the, for, relevant, services, businesses, an, expert, leading, s, and, organizational, support, internal, opportunities, are, ensuring, teams, 2, by, strong, of, to, equal, experience'
}

foreach ($row in $explanations.Keys) {
    $ws.Cells.Item($row, 4).Value = $explanations[$row]
}

